$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) status / datetime update ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-24 12:39:13"

# --- zh-cn sheet: row 3 (b.md) handoff details ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text storage instead of auto-converting to a
# real Boolean; resetting the style back to Normal drops the transient
# quote-prefix formatting so the cell ends up a plain shared string again.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-24 12:39:00"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8195e1567a62644750fe0d9803a2ffb08b5f986c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c315be35e7dc63e0db41cda7a0fd774dd588f588/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 (b.md) handoff details ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-24 12:39:13"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8195e1567a62644750fe0d9803a2ffb08b5f986c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c315be35e7dc63e0db41cda7a0fd774dd588f588/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
